# "Fixed light saturation in spotlight overlap" + "Added snowdrift"
#
# Sheet1 is the 08960 mark-scheme sheet. The edits:
#  - Row 7 ("Ground - elevated snow drifts, melts"): mark raised 0 -> 10 (DONE),
#    with a new "DONE" note in E7 (snowdrift implemented).
#  - Row 6 ("Ground - grass"): drop the stale "add tree view modes" TODO note in G6.
#  - Row 4 (Globe): add a new TODO note in G4 ("fire, lightning, seasons").
#  - Row 15 (Tree Viewing Modes): add a new TODO note in F15 ("flat shaded version
#    of leaves?").
#  - Row 16 (Lighting): add a new TODO note in G16 ("flip lights for reflection") -
#    this is the "fixed light saturation in spotlight overlap" part.
#  - The selected cell moves from F15 to F12.
# All dependent formulas (D3, D31, D32, D34, D37, D53, D54, D56, and the chart
# that plots D56) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 16: Lighting - add new note in G16 (spotlight overlap saturation fix) ---
$ws.Range("G16").Value = "flip lights for reflection"

# --- Row 4: Globe - add new note "fire, lightning, seasons" in G4 ---
$ws.Range("G4").Value = "fire, lightning, seasons"

# --- Row 6: Ground - grass - remove stale note "add tree view modes" from G6 ---
$ws.Range("G6").ClearContents()

# --- Row 7: Ground - elevated snow drifts, melts - mark DONE, snowdrift added ---
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = "DONE"
$ws.Range("E7").Font.Color = $ws.Range("E4").Font.Color
$ws.Range("G7").ClearContents()

# --- Row 15: Tree Viewing Modes - add new note in F15 ---
$ws.Range("F15").Value = "flat shaded version of leaves?"

# --- Update the active selection from F15 to F12 ---
$ws.Activate()
$ws.Range("F12").Select()

# Recalculate so every dependent cell & the embedded chart pick up the new values.
$excel.CalculateFull()
